$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force all target cells to Text format so numeric-looking strings
# (e.g. "549.41", "1.00", "0.0536") are preserved exactly as text,
# matching the original inlineStr cell values instead of being
# auto-converted to floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.384.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.459.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.82%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.67"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.71%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.459.76"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.74%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -9.89%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.72%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -8.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.05"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -9.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.901.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.82%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -9.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.315.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.462.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.11"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -8.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.07"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -9.01%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "318.45"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -7.59%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.01"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0976"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -13.62%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.583.34"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.58%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "550.46"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.06%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.47"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -10.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.27"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -10.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.67"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.14%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -8.30%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -9.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -11.78%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.83"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -12.32%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.21%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.05%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "142.52"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.31%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.76"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -8.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.36"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.18"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -10.09%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -8.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.54"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -11.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0536"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -9.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0938"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.83%  "
